# Update "paises" (countries) COVID stats workbook.
# - Refresh per-country stats for several countries (India, Pakistan, Israel,
#   Kirguistan, Uzbekistan, El Salvador, Tailandia).
# - Because the sheet is kept sorted by "Casos totales" (column B) descending,
#   a few of those updates change relative rank, so the rows in between are
#   re-ordered (Pakistan now outranks Arabia Saudita; Uzbekistan now outranks
#   Serbia/Irlanda/Marruecos).
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes
$rows = @(
    @(6,   "India",          1856754, 1423, 1231576, 586185, 0, 22, 38993),
    @(16,  "Pakistan",        280461,  432,  249397,  25065, 0, 15,  5999),
    @(17,  "Arabia Saudita",  280093,    0,  242055,  35089, 0,  0,  2949),
    @(36,  "Israel",           74903,  473,   49757,  24600, 0,  0,   546),
    @(56,  "Kirguistan",       37397,  268,   28451,   7521, 0,  5,  1425),
    @(61,  "Uzbekistan",       26550,  484,   17262,   9127, 0,  4,   161),
    @(62,  "Serbia",           26451,    0,   14047,  11806, 0,  0,   598),
    @(63,  "Irlanda",          26208,    0,   23364,   1081, 0,  0,  1763),
    @(64,  "Marruecos",        26196,    0,   18968,   6827, 0,  0,   401),
    @(73,  "El Salvador",      17843,    0,    8845,   8512, 0,  9,   486),
    @(114, "Tailandia",         3321,    1,    3142,    121, 0,  0,    58)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}

# Bump the "updated at" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 07:42"
